$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data row (shifts existing
# row 1 -> row 2, carrying its values/styles/number formats with it).
$ws.Rows("1:1").Insert()

# New header row values (A1:E1)
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Firstname"
$ws.Range("D1").Value = "Lastname"
$ws.Range("E1").Value = "Role"

# New data values for the (shifted) data row (C2:E2) - A2/B2 already
# contain the original email/password values after the row insert.
$ws.Range("C2").Value = "Anu"
$ws.Range("D2").Value = "Bhat"
$ws.Range("E2").Value = "Super User"

# Give column E an explicit width (closest achievable match).
$ws.Columns("E:E").ColumnWidth = 10

# The hyperlink that used to sit on A1 needs to move down to A2 along
# with the data. This runtime does not auto-shift the Hyperlinks
# collection on row insert, so drop and recreate it pointing at A2.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:anu.bhat@harbingergroup.com") | Out-Null

# Re-apply the Hyperlink cell style to A2 (Hyperlinks.Add re-styles the
# cell but via a freshly duplicated style record; re-applying the named
# style lets the engine reuse the original "Hyperlink" style index).
$ws.Range("A2").Style = "Hyperlink"

# Match the saved selection shown in the diff.
$ws.Range("E2").Select() | Out-Null
